# Disable "smart quotes" auto-replacement so literal straight quotes/apostrophes
# inserted below are preserved exactly as typed.
try {
    $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
} catch {
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: Insert a new "Meta description" paragraph right after the first
# paragraph ("Play Curse of the Werewolf Megaways Free | Review", Heading1)
# and before the "Unleashing the Fury..." (Heading2) paragraph.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
# Pick an insertion point strictly inside paragraph 1's own text (not on the
# paragraph-boundary position) so the new paragraph is cleanly inserted
# between paragraph 1 and paragraph 2 without disturbing either of them or
# inheriting their styles/rsids.
$midPos = $p1.Range.Start + 5

$metaXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Curse of the Werewolf Megaways, a high volatility slot game with multiple bonus features and up to 46,656 ways to win. Play for free today.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertTarget = $d.Range($midPos, $midPos)
[void]$insertTarget.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# Part 2: Near the end of the document, remove the duplicate bold
# "Play Curse of the Werewolf Megaways Free | Review" paragraph, and replace
# the text of the following italic paragraph with the new image prompt text.
# Locate the two paragraphs by content (searching from the end) rather than
# by a hard-coded index, so the logic is resilient to the exact paragraph
# count.
# ---------------------------------------------------------------------------
$boldDupeIndex = -1
$italicIndex = -1
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($boldDupeIndex -eq -1 -and $text -like "*Play Curse of the Werewolf Megaways Free | Review*") {
        $boldDupeIndex = $i
    }
    if ($italicIndex -eq -1 -and $text -like "*Read our review of Curse of the Werewolf Megaways, a high volatility slot game*") {
        $italicIndex = $i
    }
    if ($boldDupeIndex -ne -1 -and $italicIndex -ne -1) {
        break
    }
}

$boldDupe = $d.Paragraphs.Item($boldDupeIndex)
$boldDupe.Range.Delete()

# Paragraph indices after the above paragraph shift down by one once it is
# removed (only relevant if the bold paragraph preceded the italic one).
if ($boldDupeIndex -lt $italicIndex) {
    $italicIndex = $italicIndex - 1
}

$italicPara = $d.Paragraphs.Item($italicIndex)
$italicRange = $italicPara.Range

$oldText = "Read our review of Curse of the Werewolf Megaways, a high volatility slot game with multiple bonus features and up to 46,656 ways to win. Play for free today."
$found = $italicRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$newText = 'Prompt: Create a cartoon-style feature image for "Curse of the Werewolf Megaways" featuring a happy Maya warrior with glasses. For this feature image, DALLE could create a cartoon-style illustration that captures the mood and theme of "Curse of the Werewolf Megaways." The image could feature a happy Maya warrior with glasses, adding a fun and unexpected twist to the werewolf legend. The Maya warrior could be depicted in vibrant colors, wearing traditional garments and accessories like feathered headdresses and intricate jewelry. They could be shown in a dynamic pose, holding a weapon or shield, as though they are ready to take on the werewolf or protect their village. The werewolf could also be included in the image, perhaps lurking in the background or partially hidden behind a building or tree. The Maya warrior could be shown confidently facing the werewolf, unafraid and ready to take on the challenge. Overall, the image should be attention-grabbing and visually appealing, incorporating elements of the game''s theme and characters in a creative and playful way.'

$italicRange.Text = $newText
